# Weekly update: a brand-new week of "Perejil" price data is added at the
# top of the table (rows 8-9), every subsequent weekly pair of rows
# (Primera/Segunda) is pushed down into the slot of the following pair,
# and the pair that falls off the bottom is appended as two new rows
# (156-157) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 8
$lastRow  = 155
$numPairs = ($lastRow - $startRow + 1) / 2   # 74 weekly pairs (Primera/Segunda)

# Columns that carry the weekly data and therefore shift down one pair:
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg).
$shiftCols = @(4, 10, 11, 12, 13, 14, 15, 16)

# --- 1. Snapshot all current ("old") values for every pair, per column ----
$old = @{}
foreach ($col in $shiftCols) {
    $row0vals = @()
    $row1vals = @()
    for ($i = 0; $i -lt $numPairs; $i++) {
        $r0 = $startRow + 2 * $i
        $r1 = $r0 + 1
        $row0vals += ,$ws.Cells.Item($r0, $col).Value2
        $row1vals += ,$ws.Cells.Item($r1, $col).Value2
    }
    $old[$col] = @($row0vals, $row1vals)
}

# --- 2. Append two brand-new rows (156-157) for the pair that falls off --
#        the bottom (the last pair's OLD data). Copy row 154-155 first so
#        styles/number-formats/borders are carried over, then the loop
#        below will overwrite the shifting-column values as needed.
$newRow0 = $lastRow + 1
$newRow1 = $lastRow + 2
$ws.Range(("A" + ($lastRow - 1) + ":R" + $lastRow)).Copy($ws.Range("A" + $newRow0 + ":R" + $newRow1))

# --- 3. Write the shifted values back, pair by pair -----------------------
#    New pair i (i = 1..numPairs-1) gets OLD pair (i-1)'s values.
#    New pair numPairs (the appended rows) gets OLD pair (numPairs-1)'s values.
#    New pair 0 (rows 8-9) gets brand-new data (see step 4).
for ($i = 1; $i -le $numPairs; $i++) {
    if ($i -lt $numPairs) {
        $r0 = $startRow + 2 * $i
        $r1 = $r0 + 1
    } else {
        $r0 = $newRow0
        $r1 = $newRow1
    }
    $srcIdx = $i - 1
    foreach ($col in $shiftCols) {
        $ws.Cells.Item($r0, $col).Value = $old[$col][0][$srcIdx]
        $ws.Cells.Item($r1, $col).Value = $old[$col][1][$srcIdx]
    }
}

# --- 4. New pair 0 (rows 8-9): a brand-new week's worth of data -----------
#    Volumen (J), Unidad (N) and Origen (O) are unchanged from what they
#    were; only Fecha (D) and the three price columns (K,L,M) plus the
#    derived Precio $/Kg (P) change.
$ws.Cells.Item($startRow, 4).Value     = 44756
$ws.Cells.Item($startRow, 11).Value    = 700
$ws.Cells.Item($startRow, 12).Value    = 800
$ws.Cells.Item($startRow, 13).Value    = 750
$ws.Cells.Item($startRow, 16).Value    = 750

$ws.Cells.Item($startRow + 1, 4).Value  = 44756
$ws.Cells.Item($startRow + 1, 11).Value = 600
$ws.Cells.Item($startRow + 1, 12).Value = 600
$ws.Cells.Item($startRow + 1, 13).Value = 600
$ws.Cells.Item($startRow + 1, 16).Value = 600

Write-Host "Done shifting $numPairs pairs and appending rows $newRow0-$newRow1"
